# Derek's Log - add the week-of-2017-01-09 (Monday) entries to the "Logs" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- Row 739: new "MONDAY" divider row. Clone the existing divider row (737)
#     to pick up its fill/border styling, then swap the day-name text. ---
$ws.Range("A737:F737").Copy($ws.Range("A739:F739"))
$ws.Range("C739").Value = "MONDAY"

# --- Rows 740-744, 748-752: normal data rows (A..F all styled).
#     Row 736 has that exact style pattern (10/21/17/18/19/20), clone it. ---
$ws.Range("A736:F736").Copy($ws.Range("A740:F740"))
$ws.Range("A736:F736").Copy($ws.Range("A741:F741"))
$ws.Range("A736:F736").Copy($ws.Range("A742:F742"))
$ws.Range("A736:F736").Copy($ws.Range("A743:F743"))
$ws.Range("A736:F736").Copy($ws.Range("A744:F744"))
$ws.Range("A736:F736").Copy($ws.Range("A748:F748"))
$ws.Range("A736:F736").Copy($ws.Range("A749:F749"))
$ws.Range("A736:F736").Copy($ws.Range("A750:F750"))
$ws.Range("A736:F736").Copy($ws.Range("A751:F751"))
$ws.Range("A736:F736").Copy($ws.Range("A752:F752"))

# --- Row 745: identical in every column to row 738 except the date -
#     clone row 738 (keeps the rich-text "Door code 11012* " comment intact). ---
$ws.Range("A738:F738").Copy($ws.Range("A745:F745"))

# --- Rows 746-747: data rows with no comment (F left empty).
#     Row 734 has that exact style pattern (10/21/17/18/19, no F), clone it. ---
$ws.Range("A734:E734").Copy($ws.Range("A746:E746"))
$ws.Range("A734:E734").Copy($ws.Range("A747:E747"))

# --- Fill in the actual values for each new row ---

$ws.Range("A740").Value = "Demo"
$ws.Range("B740").Value = 42744
$ws.Range("C740").Value = "1550"
$ws.Range("D740").Value = "VH"
$ws.Range("E740").Value = "1158"
$ws.Range("F740").Value = "Meet instructor Sanjeev Dhuga"

$ws.Range("A741").Value = "Demo"
$ws.Range("B741").Value = 42744
$ws.Range("C741").Value = "1550"
$ws.Range("D741").Value = "VH"
$ws.Range("E741").Value = "D"
$ws.Range("F741").Value = "Meet instructor Douglas McCready"

$ws.Range("A742").Value = "AV Shutdown"
$ws.Range("B742").Value = 42744
$ws.Range("C742").Value = "1630"
$ws.Range("D742").Value = "MC"
$ws.Range("E742").Value = "101A"
$ws.Range("F742").Value = "Pick up wireless keyboard and TV remote control. To FDRS 164."

$ws.Range("A743").Value = "Demo"
$ws.Range("B743").Value = 42744
$ws.Range("C743").Value = "1650"
$ws.Range("D743").Value = "MC"
$ws.Range("E743").Value = "140"
$ws.Range("F743").Value = "PC and neck mic is there.  Demo for Khan Le. Door code 7083*"

$ws.Range("A744").Value = "Demo"
$ws.Range("B744").Value = 42744
$ws.Range("C744").Value = "1720"
$ws.Range("D744").Value = "R"
$ws.Range("E744").Value = "S205"
$ws.Range("F744").Value = "Meet instructor Susan Ehrlich"

# Row 745 only needs the date changed - everything else already matches row 738.
$ws.Range("B745").Value = 42744

$ws.Range("A746").Value = "AV Shutdown"
$ws.Range("B746").Value = 42744
$ws.Range("C746").Value = "1730"
$ws.Range("D746").Value = "R"
$ws.Range("E746").Value = "N203"

$ws.Range("A747").Value = "AV Shutdown"
$ws.Range("B747").Value = 42744
$ws.Range("C747").Value = "1830"
$ws.Range("D747").Value = "R"
$ws.Range("E747").Value = "S203"

$ws.Range("A748").Value = "AV Shutdown"
$ws.Range("B748").Value = 42744
$ws.Range("C748").Value = "1830"
$ws.Range("D748").Value = "R"
$ws.Range("E748").Value = "N102"
$ws.Range("F748").Value = "Nat Taylor Cinema. Lock all cinema doors after shutdown."

$ws.Range("A749").Value = "Demo"
$ws.Range("B749").Value = 42744
$ws.Range("C749").Value = "1850"
$ws.Range("D749").Value = "R"
$ws.Range("E749").Value = "S103"
$ws.Range("F749").Value = "Meet instructor Karen Murray"

$ws.Range("A750").Value = "Demo"
$ws.Range("B750").Value = 42744
$ws.Range("C750").Value = "1850"
$ws.Range("D750").Value = "SLH"
$ws.Range("E750").Value = "A"
$ws.Range("F750").Value = "Meet instructor Jill Prindiville"

$ws.Range("A751").Value = "Pickup PC"
$ws.Range("B751").Value = 42744
$ws.Range("C751").Value = "1930"
$ws.Range("D751").Value = "MC"
$ws.Range("E751").Value = "140"
$ws.Range("F751").Value = "Leave portable screen. Door code 7083* return to FDRS 156A."

$ws.Range("A752").Value = "Pickup Mic"
$ws.Range("B752").Value = 42744
$ws.Range("C752").Value = "1930"
$ws.Range("D752").Value = "MC"
$ws.Range("E752").Value = "140"
$ws.Range("F752").Value = "Neck mic and small PA to FDRS 156A."

# --- Update the sheet's last selection to match where Excel would have
#     left the cursor after typing in the new rows. ---
$ws.Range("A752").Select()
